# Add version numbers to schemas.
#
# 1. Insert two new leading columns ("version", "description") into the
#    "Export as TSV" sheet, shifting every existing column (and its
#    comment / data validation) two places to the right.
# 2. Add a new "version list" worksheet (right after "Export as TSV")
#    that supplies the allowed value ("1") for the new "version" column.
# 3. Wire up a data validation on the new "version" column that points
#    at the new "version list" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export as TSV")

# ---------------------------------------------------------------------
# Step 0: remember the existing header comments (keyed by their current
# 1-based column number) so we can re-create them, shifted, after the
# column insert below. (Comments do not automatically move when columns
# are inserted.)
# ---------------------------------------------------------------------
$oldComments = @{}
$lastCol = $ws.UsedRange.Columns.Count
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cmt = $cell.Comment
    if ($cmt -ne $null) {
        $oldComments[$c] = $cmt.Text()
        $cmt.Delete()
    }
}

# ---------------------------------------------------------------------
# Step 1: insert the two new columns at the front of the sheet.
# ---------------------------------------------------------------------
$ws.Range("A1:B1").EntireColumn.Insert()

# Match the existing bold / centered / wrapped header formatting.
$ws.Range("A1:B1").Font.Bold = $true
$ws.Range("A1:B1").HorizontalAlignment = -4108
$ws.Range("A1:B1").WrapText = $true

$ws.Range("A1").Value = "version"
$ws.Range("B1").Value = "description"

# ---------------------------------------------------------------------
# Step 2: re-create the old comments, shifted two columns to the right,
# then add the two new comments for the "version" / "description"
# columns.
# ---------------------------------------------------------------------
foreach ($c in $oldComments.Keys) {
    $newCol = $c + 2
    $cell = $ws.Cells.Item(1, $newCol)
    $cell.AddComment($oldComments[$c]) | Out-Null
}

$ws.Cells.Item(1, 1).AddComment("Version of the schema to use when validating this metadata.") | Out-Null
$ws.Cells.Item(1, 2).AddComment("Free-text description of this assay.") | Out-Null

# ---------------------------------------------------------------------
# Step 3: add the new "version list" worksheet right after
# "Export as TSV", with the single allowed value "1".
# ---------------------------------------------------------------------
$versionSheet = $wb.Worksheets.Add($null, $ws)
$versionSheet.Name = "version list"
$versionSheet.Range("A1").NumberFormat = "@"
$versionSheet.Range("A1").Value = "1"

# ---------------------------------------------------------------------
# Step 4: add the data validation for the new "version" column.
# ---------------------------------------------------------------------
$versionValidation = $ws.Range("A2:A1048576").Validation
$versionValidation.Add(3, 1, 1, "='version list'!`$A`$1:`$A`$1")
$versionValidation.ErrorTitle = "Value must come from list"
$versionValidation.ErrorMessage = "Value must be one of: 1."
$versionValidation.ShowInput = $true
$versionValidation.ShowError = $true
